# SistemasOrdenanzas.xlsx - add the VALVERDE municipality vehicle-tax table
# (new ordinance rows appended below the existing BENAVIDES / ASTORGA / TRUCHAS blocks).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row, TIPOVEHICULO, UNIDAD, MINIMO, MAXIMO, IMPORTE  (AYUNTAMIENTO column is always VALVERDE)
$rows = @(
    @(103, "TURISMO",     "CABALLOS", 0,       7.99,    20),
    @(104, "TURISMO",     "CABALLOS", 8,       11.99,   30),
    @(105, "TURISMO",     "CABALLOS", 12,      15.99,   40),
    @(106, "TURISMO",     "CABALLOS", 16,      19.99,   50),
    @(107, "TURISMO",     "CABALLOS", 20,      9999,    60),

    @(110, "AUTOBUS",     "PLAZAS",   0,       20,      80),
    @(111, "AUTOBUS",     "PLAZAS",   21,      50,      100),
    @(112, "AUTOBUS",     "PLAZAS",   51,      9999,    120),

    @(114, "CAMION",      "KG",       0,       999,     45),
    @(115, "CAMION",      "KG",       1000,    2999,    60),
    @(116, "CAMION",      "KG",       3000,    9999,    95),
    @(117, "CAMION",      "KG",       10000,   9999999, 120),

    @(119, "TRACTOR",     "CABALLOS", 0,       15.99,   20),
    @(120, "TRACTOR",     "CABALLOS", 16,      25,      30),
    @(121, "TRACTOR",     "CABALLOS", 25.01,   9999,    72),

    @(123, "REMOLQUE",    "KG",       751,     999,     37),
    @(124, "REMOLQUE",    "KG",       1000,    2999,    45),
    @(125, "REMOLQUE",    "KG",       3000,    9999999, 91),

    @(127, "CICLOMOTOR",  "CC",       0,       49.99,   12),

    @(128, "MOTOCICLETA", "CC",       0,       125,     20),
    @(129, "MOTOCICLETA", "CC",       125.01,  250,     30),
    @(130, "MOTOCICLETA", "CC",       250.01,  500,     45),
    @(131, "MOTOCICLETA", "CC",       500.01,  1000,    56),
    @(132, "MOTOCICLETA", "CC",       1000.01, 9999,    67)
)

foreach ($r in $rows) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value = "VALVERDE"
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
}

# Match the final selection from the authored edit (Excel auto-scrolls the
# view to keep the active cell visible when the sheet is saved).
[void]$ws.Range("F134").Select()
